$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: Conversion del dia text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$lines = @(
    "Conversión del día 💰",
    "✅ Dólar paralelo: 68",
    "",
    "Binance",
    "✅ 1000 Bs = 12.05 = 48327.71 pesos",
    "✅ 48327.71 pesos = 12.03 = 965.11 Bs",
    "",
    "Promedio competencia",
    "✅ Tasa pesos: 20",
    "✅ Tasa Bs: 20",
    "✅ % Ganancia: 20%"
)
$newText = [string]::Join("`n", $lines)

$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet: N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 83
$wsTasas.Range("O10").Value = 4011.2
$wsTasas.Range("N12").Value = 4016
$wsTasas.Range("O12").Value = 80.2
